$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column C (td_sim_1) and column D (record_atd) values for rows 2-64
# per corrected relevance markers Appenzeller-Herzog (2019) - van Dis (2020)

$ws.Cells.Item(2, 3).Value = 395
$ws.Cells.Item(2, 4).Value = 434
$ws.Cells.Item(3, 3).Value = 12
$ws.Cells.Item(3, 4).Value = 10
$ws.Cells.Item(4, 3).Value = 92
$ws.Cells.Item(4, 4).Value = 79.5
$ws.Cells.Item(5, 3).Value = 873
$ws.Cells.Item(5, 4).Value = 1029.5
$ws.Cells.Item(6, 3).Value = 143
$ws.Cells.Item(6, 4).Value = 160
$ws.Cells.Item(7, 3).Value = 51
$ws.Cells.Item(7, 4).Value = 58
$ws.Cells.Item(8, 3).Value = 22
$ws.Cells.Item(8, 4).Value = 18.5
$ws.Cells.Item(9, 3).Value = 135
$ws.Cells.Item(9, 4).Value = 159.5
$ws.Cells.Item(10, 3).Value = 1555
$ws.Cells.Item(10, 4).Value = 1496.5
$ws.Cells.Item(11, 3).Value = 540
$ws.Cells.Item(11, 4).Value = 532
$ws.Cells.Item(12, 3).Value = 150
$ws.Cells.Item(12, 4).Value = 151.5
$ws.Cells.Item(13, 3).Value = 142
$ws.Cells.Item(13, 4).Value = 163.5
$ws.Cells.Item(14, 3).Value = 271
$ws.Cells.Item(14, 4).Value = 269
$ws.Cells.Item(15, 3).Value = 27
$ws.Cells.Item(15, 4).Value = 23
$ws.Cells.Item(16, 3).Value = 16
$ws.Cells.Item(16, 4).Value = 14
$ws.Cells.Item(17, 3).Value = 56
$ws.Cells.Item(17, 4).Value = 62
$ws.Cells.Item(18, 3).Value = 187
$ws.Cells.Item(18, 4).Value = 199
$ws.Cells.Item(19, 3).Value = 3011
$ws.Cells.Item(19, 4).Value = 2907
$ws.Cells.Item(20, 3).Value = 2281
$ws.Cells.Item(20, 4).Value = 2672
$ws.Cells.Item(21, 3).Value = 909
$ws.Cells.Item(21, 4).Value = 984.5
$ws.Cells.Item(22, 3).Value = 11
$ws.Cells.Item(22, 4).Value = 10
$ws.Cells.Item(23, 3).Value = 7
$ws.Cells.Item(23, 4).Value = 8.5
$ws.Cells.Item(24, 3).Value = 203
$ws.Cells.Item(24, 4).Value = 215
$ws.Cells.Item(25, 3).Value = 58
$ws.Cells.Item(25, 4).Value = 65.5
$ws.Cells.Item(26, 3).Value = 189
$ws.Cells.Item(26, 4).Value = 202.5
$ws.Cells.Item(27, 3).Value = 121
$ws.Cells.Item(27, 4).Value = 135
$ws.Cells.Item(28, 3).Value = 61
$ws.Cells.Item(28, 4).Value = 68
$ws.Cells.Item(29, 3).Value = 487
$ws.Cells.Item(29, 4).Value = 563
$ws.Cells.Item(30, 3).Value = 218
$ws.Cells.Item(30, 4).Value = 211.5
$ws.Cells.Item(31, 3).Value = 199
$ws.Cells.Item(31, 4).Value = 206.5
$ws.Cells.Item(32, 3).Value = 50
$ws.Cells.Item(32, 4).Value = 53
$ws.Cells.Item(33, 3).Value = 54
$ws.Cells.Item(33, 4).Value = 61.5
$ws.Cells.Item(34, 3).Value = 19
$ws.Cells.Item(34, 4).Value = 16
$ws.Cells.Item(35, 3).Value = 615
$ws.Cells.Item(35, 4).Value = 611.5
$ws.Cells.Item(36, 3).Value = 687
$ws.Cells.Item(36, 4).Value = 629
$ws.Cells.Item(37, 3).Value = 3
$ws.Cells.Item(37, 4).Value = 25
$ws.Cells.Item(38, 3).Value = 107
$ws.Cells.Item(38, 4).Value = 108
$ws.Cells.Item(39, 3).Value = 375
$ws.Cells.Item(39, 4).Value = 357.5
$ws.Cells.Item(40, 3).Value = 1509
$ws.Cells.Item(40, 4).Value = 1479.5
$ws.Cells.Item(41, 3).Value = 447
$ws.Cells.Item(41, 4).Value = 466.5
$ws.Cells.Item(42, 3).Value = 41
$ws.Cells.Item(42, 4).Value = 34
$ws.Cells.Item(43, 3).Value = 378
$ws.Cells.Item(43, 4).Value = 378
$ws.Cells.Item(44, 3).Value = 23
$ws.Cells.Item(44, 4).Value = 18.5
$ws.Cells.Item(45, 3).Value = 1924
$ws.Cells.Item(45, 4).Value = 1748
$ws.Cells.Item(46, 3).Value = 155
$ws.Cells.Item(46, 4).Value = 148.5
$ws.Cells.Item(47, 3).Value = 202
$ws.Cells.Item(47, 4).Value = 215
$ws.Cells.Item(48, 3).Value = 85
$ws.Cells.Item(48, 4).Value = 87
$ws.Cells.Item(49, 3).Value = 186
$ws.Cells.Item(49, 4).Value = 194
$ws.Cells.Item(50, 3).Value = 400
$ws.Cells.Item(50, 4).Value = 397.5
$ws.Cells.Item(51, 3).Value = 193
$ws.Cells.Item(51, 4).Value = 204
$ws.Cells.Item(52, 3).Value = 40
$ws.Cells.Item(52, 4).Value = 37
$ws.Cells.Item(53, 3).Value = 42
$ws.Cells.Item(53, 4).Value = 48
$ws.Cells.Item(54, 3).Value = 84
$ws.Cells.Item(54, 4).Value = 97
$ws.Cells.Item(55, 3).Value = 205
$ws.Cells.Item(55, 4).Value = 215.5
$ws.Cells.Item(56, 3).Value = 8
$ws.Cells.Item(56, 4).Value = 23
$ws.Cells.Item(57, 3).Value = 228
$ws.Cells.Item(57, 4).Value = 236
$ws.Cells.Item(58, 3).Value = 55
$ws.Cells.Item(58, 4).Value = 67
$ws.Cells.Item(59, 3).Value = 72
$ws.Cells.Item(59, 4).Value = 83
$ws.Cells.Item(60, 3).Value = 53
$ws.Cells.Item(60, 4).Value = 59.5
$ws.Cells.Item(61, 3).Value = 28
$ws.Cells.Item(61, 4).Value = 25.5
$ws.Cells.Item(62, 3).Value = 15
$ws.Cells.Item(62, 4).Value = 11
$ws.Cells.Item(63, 3).Value = 43
$ws.Cells.Item(63, 4).Value = 53
$ws.Cells.Item(64, 3).Value = 334.6451612903226
